# The commit adds one new weekly price-report row for Palta (Hass, Primera)
# right before the existing row 577, pushing every following row down by
# one (old row 577 -> new row 578, ..., old row 677 -> new row 678).
#
# Resolve the workbook/sheet via $excel (note: $wb.ActiveSheet is not
# reliable in this runtime, but $excel.ActiveWorkbook is).
$wb2 = $excel.ActiveWorkbook
$ws = $wb2.Worksheets.Item(1)

# Insert a fresh row above the current row 577; this shifts rows
# 577..677 down to 578..678 and keeps their formatting/values intact.
$ws.Range("A577").EntireRow.Insert()

# Populate the newly inserted row 577 with the new record.
$ws.Range("A577").Value = 7
$ws.Range("B577").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C577").Value = 'Ñuble'
$ws.Range("D577").Value = 44951
$ws.Range("E577").Value = 16
$ws.Range("F577").Value = 'Fruta'
$ws.Range("G577").Value = 100106
$ws.Range("H577").Value = 'Oleaginosos'
$ws.Range("I577").Value = 100106002
$ws.Range("J577").Value = 'Palta'
$ws.Range("K577").Value = 'Hass'
$ws.Range("L577").Value = 'Primera'
$ws.Range("M577").Value = 400
$ws.Range("N577").Value = 3400
$ws.Range("O577").Value = 3500
$ws.Range("P577").Value = 3450
$ws.Range("Q577").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R577").Value = 'Provincia de Quillota'
$ws.Range("S577").Value = 3450
$ws.Range("T577").Value = 1
